$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet ships with cell-edit protection enabled; temporarily lift it so
# the Percent Change values can be written, then restore protection.
$ws.Unprotect()

$values = @(
    0.008836896145430151,
    0.01486011057128489,
    0.002250574530469951,
    0.005669606664236326,
    0.01795521757498908,
    -0.003101309441764344,
    -0.007986024457199981,
    0.005285013212533052,
    0.008166295471417895,
    0.008844953173777315,
    -0.006491297692612807,
    0.0132180667794275,
    -0.00231770369049733,
    0.008572959457045837,
    0.00117332638028822,
    -0.005849197064766498,
    0.008175466397752684,
    0.0003122853038535212,
    -0.01578616807981259,
    0.03623013962445842,
    0.015343083652269,
    0.00152710613387641,
    0.007501103103397577,
    -0.007232767232767334,
    0.003412470140886414,
    0.002127282396738206,
    0.001669449081802998,
    -0.005008923944959442,
    0.004739105398772203
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $values[$i]
}

$ws.Protect()
